# Apply the edits described by the commit diff:
#  1. Update the product name in cell A3 from "Smart TV" to "Mobiles 5g"
#     (this changes the shared string table entry used by A3).
#  2. Move the active cell selection on the sheet from A7 to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Mobiles 5g"

$ws.Range("A6").Select() | Out-Null
